$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A17").Value = "edit1"
$ws.Range("B17").Value = "riya-morankar"
$ws.Range("C17").Value = "Merged"
$ws.Range("D17").Value = "N/A"

# "2025-06-19" looks like a real date, so Excel would normally auto-convert
# it to a date serial when assigned to a General-formatted cell. Force the
# cell to Text first so the literal string is preserved, then restore the
# default "Normal" style so no extra formatting is left behind.
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "2025-06-19"
$ws.Range("E17").Style = "Normal"

$ws.Range("F17").Value = "37dcd755ab7f59fccebf0f22a88d047d73a1c753"
